$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EVCRSbRIC")

# Insert a new column before column L, splitting the merged
# "ISIC 20T21" header into two separate ISIC codes: "ISIC 20" (col K)
# and "ISIC 21" (new col L). All existing data from column L onward
# shifts one column to the right.
$ws.Range("L:L").Insert()

# K1 previously held "ISIC 20T21" -> now just "ISIC 20"
$ws.Range("K1").Value = "ISIC 20"

# The newly inserted L1 gets the other half of the split header
$ws.Range("L1").Value = "ISIC 21"

# The newly inserted L2 data cell mirrors K2 (both were 0 under the
# old combined "ISIC 20T21" column)
$ws.Range("L2").Value = 0
